$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.513.95"
$ws.Range("E2").Value = "  +1.48%  "

$ws.Range("D3").Value = "1.826.12"
$ws.Range("E3").Value = "  +1.83%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'317.03"
$ws.Range("E5").Value = "  -0.01%  "

$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("D7").Value = "'0.5410"
$ws.Range("E7").Value = "  +0.23%  "

$ws.Range("D8").Value = "'0.4048"
$ws.Range("E8").Value = "  +7.34%  "

$ws.Range("D9").Value = "'0.07665"
$ws.Range("E9").Value = "  +3.16%  "

$ws.Range("E10").Value = "  +2.44%  "

$ws.Range("D11").Value = "'41.87"
$ws.Range("E11").Value = "  +0.44%  "

$ws.Range("D12").Value = "'6.326"
$ws.Range("E12").Value = "  +3.64%  "

$ws.Range("D13").Value = "'7.638"
$ws.Range("E13").Value = "  +5.78%  "

$ws.Range("E14").Value = "  +2.01%  "

$ws.Range("E15").Value = "  -0.05%  "

$ws.Range("D16").Value = "1.825.79"
$ws.Range("E16").Value = "  +2.06%  "

$ws.Range("E17").Value = "  +3.04%  "

$ws.Range("D18").Value = "'89.96"
$ws.Range("E18").Value = "  +1.02%  "

$ws.Range("D19").Value = "'0.06602"
$ws.Range("E19").Value = "  +1.98%  "

$ws.Range("E20").Value = "  +2.59%  "

$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("D22").Value = "'6.071"
$ws.Range("E22").Value = "  +2.94%  "

$ws.Range("D23").Value = "28.524.11"
$ws.Range("E23").Value = "  +1.48%  "

$ws.Range("D24").Value = "'11.17"
$ws.Range("E24").Value = "  -0.04%  "

$ws.Range("D25").Value = "'2.272"
$ws.Range("E25").Value = "  +8.58%  "

$ws.Range("D26").Value = "'157.84"
$ws.Range("E26").Value = "  +1.87%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'20.74"
$ws.Range("E27").Value = "  +2.45%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.448"
$ws.Range("E28").Value = "  +7.10%  "

$ws.Range("D29").Value = "2.034.81"
$ws.Range("E29").Value = "  +2.08%  "

$ws.Range("D30").Value = "'123.93"
$ws.Range("E30").Value = "  +2.70%  "

$ws.Range("E31").Value = "  +0.34%  "

$ws.Range("E32").Value = "  +4.80%  "

$ws.Range("D33").Value = "'5.678"
$ws.Range("E33").Value = "  +2.37%  "

$ws.Range("D34").Value = "'0.07432"
$ws.Range("E34").Value = "  +14.32%  "

$ws.Range("E35").Value = "  -0.42%  "

$ws.Range("D36").Value = "'0.2240"
$ws.Range("E36").Value = "  -0.74%  "

$ws.Range("D37").Value = "'0.02346"
$ws.Range("E37").Value = "  +2.59%  "

$ws.Range("E38").Value = "  +4.07%  "

$ws.Range("D39").Value = "'8.887"
$ws.Range("E39").Value = "  +5.04%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.6296"
$ws.Range("E40").Value = "  +2.23%  "

$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'11.36"
$ws.Range("E41").Value = "  +2.69%  "

$ws.Range("D42").Value = "'1.189"
$ws.Range("E42").Value = "  +1.49%  "

$ws.Range("D43").Value = "'0.9996"
$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("E44").Value = "  -3.50%  "

$ws.Range("D45").Value = "'13.45"
$ws.Range("E45").Value = "  +1.66%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5873"
$ws.Range("E46").Value = "  +1.64%  "

$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.699"
$ws.Range("E47").Value = "  +0.72%  "

$ws.Range("D48").Value = "'125.51"
$ws.Range("E48").Value = "  -0.01%  "

$ws.Range("D49").Value = "'2.004"
$ws.Range("E49").Value = "  +4.50%  "

$ws.Range("D50").Value = "'1.200"
$ws.Range("E50").Value = "  +1.01%  "

$ws.Range("D51").Value = "'0.06889"
$ws.Range("E51").Value = "  +1.16%  "
